$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '39.511.72'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.00%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.165.27'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.07%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.633'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.68%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.64'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.01%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("E9").Value = '  +1.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0852'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.44%  '

$ws.Range("E11").Value = '  +0.46%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.08'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.13%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.488.79'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.811'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.12%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.51'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.24%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.140.54'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.24%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '39.508.99'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.94%  '

$ws.Range("E19").Value = '  +1.45%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.45%  '

$ws.Range("E21").Value = '  +0.93%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '230.17'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.91%  '

$ws.Range("E23").Value = '  +0.07%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.33'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.11%  '

$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.33%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.70'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.62%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '172.08'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.00%  '

$ws.Range("E28").Value = '  +0.60%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.88'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.78%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.43'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.70%  '

$ws.Range("E31").Value = '  +5.57%  '

$ws.Range("E32").Value = '  +2.09%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.60'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.45%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.71'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.65%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.04'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.23%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0620'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.07%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.69'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.01%  '

$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.43'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.67%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.17%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.95'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.95%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '102.63'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0228'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.23%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.523.76'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.43%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.21'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.59%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.43'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.25%  '

$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.11'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.31%  '

$ws.Range("B47").Value = 'HuobiToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.82'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.33%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0925'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.43%  '

$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.75'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.65%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.372.01'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.21%  '

$ws.Range("E51").Value = '  -0.54%  '
